$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Internal.Public"
$ws.Range("D3").Value = "External.Public"
